$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "https://intellisoft-consulting.github.io/igs/ChanjoKe-FHIR-IG/CodeSystem/KHTS.A"
$meta.Range("B5").Value = "KenyaCounties CodeSystem for Counties in Kenya"
$meta.Range("B8").Value = "2024-08-27T20:30:12+00:00"
$meta.Range("B9").Value = "Intellisoft Consulting Ltd"
$meta.Range("B10").Value = "Intellisoft Consulting Ltd (https://www.intellisoftkenya.com/, info[at]intellisoftkenya.com)"
$meta.Range("B12").Value = "CodeSystem for Counties in Kenya"

# --- Concepts sheet updates ---
$concepts = $wb.Worksheets.Item("Concepts")

$concepts.Range("B2").Value = "C-001"
$concepts.Range("B3").Value = "C-002"
$concepts.Range("B4").Value = "C-003"
$concepts.Range("B5").Value = "C-004"
